$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stash clean format templates (col A / col B / col C cell styles) in a safe, far-away area ---
$ws.Range("A3").Copy()
$ws.Range("AA1000").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("AA1001").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("AA1002").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Clear only the part of the table that actually changes (rows 10-23); rows 1-9 stay untouched ---
$ws.Range("A10:C23").Clear()
$ws.Rows("10:23").AutoFit()

# --- Repopulate rows 10-25 in final top-to-bottom, left-to-right order so the shared-string table ---
# --- gets rebuilt in exactly the order the target workbook uses ---
$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer os conceitos fundamentais sobre dispositivos semicondutores aplicados em circuitos eletrônicos, que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse envolvendo os principais componentes eletrônicos, como diodos, transistores e amplificadores operacionais, com escolhas adequadas de hipóteses e aplicação de ferramentas correspondentes de solução; (c) introduzir os componentes, técnicas, softwares e equipamentos utilizados na análise e projeto de circuitos eletrônicos; e (d) aplicar e estender os conceitos físicos aprendidos previamente.'
$ws.Range("C10").Value = 'Como parte fundamental da formação específica e geral, a disciplina tem por objetivos (a) fornecer os conceitos fundamentais sobre dispositivos semicondutores aplicados em circuitos eletrônicos, que são importantes para a formação em engenharia física; (b) capacitar o aluno, trabalhando individualmente e em grupo, a modelar e resolver problemas de interesse envolvendo os principais componentes eletrônicos, como diodos, transistores e amplificadores operacionais, com escolhas adequadas de hipóteses e aplicação de ferramentas correspondentes de solução; (c) introduzir os componentes, técnicas, softwares e equipamentos utilizados na análise e projeto de circuitos eletrônicos; e (d) aplicar e estender os conceitos físicos aprendidos previamente.'
$ws.Rows(10).RowHeight = 60

$ws.Range("A11").Value = 'Objectives:'
$ws.Range("B11").Value = 'As an important requirement for the specific and the general engineering studies, the course aims to (a) provide the fundamental concepts concerning the applications of semiconductor devices within the context of electronics, which are important for the learning process in physics engineering; (b) enable the student, to work as individually as in groups, to model and solve problems of interest regarding the most important electronic devices, such as diodes, transistors, and operational amplifiers, including the adequate choice of hypotheses and the application of suitable solution tools; (c) to introduce the components, techniques, software, and equipment employed in the analysis and design of electrical circuits; and (d) apply and extend the previously learned physical concepts.'
$ws.Range("C11").Value = 'As an important requirement for the specific and the general engineering studies, the course aims to (a) provide the fundamental concepts concerning the applications of semiconductor devices within the context of electronics, which are important for the learning process in physics engineering; (b) enable the student, to work as individually as in groups, to model and solve problems of interest regarding the most important electronic devices, such as diodes, transistors, and operational amplifiers, including the adequate choice of hypotheses and the application of suitable solution tools; (c) to introduce the components, techniques, software, and equipment employed in the analysis and design of electrical circuits; and (d) apply and extend the previously learned physical concepts.'
$ws.Rows(11).RowHeight = 60

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("B13").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C13").Value = '519033 - Carlos Yujiro Shigue'

$ws.Range("B14").Value = '7290967 - Emerson Gonçalves de Melo'
$ws.Range("C14").Value = '7290967 - Emerson Gonçalves de Melo'

$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = 'Introdução aos materiais e dispositivos semicondutores; Diodo; Transistor bipolar de junção; Transistores de efeito de campo; Amplificadores operacionais; Amplificadores de pequenos sinais; Fontes de alimentação. Análise e projeto de circuitos eletrônicos utilizando softwares EDA. Análises experimentais de circuitos eletrônicos.'
$ws.Range("C15").Value = 'Introdução aos materiais e dispositivos semicondutores; Diodo; Transistor bipolar de junção; Transistores de efeito de campo; Amplificadores operacionais; Amplificadores de pequenos sinais; Fontes de alimentação. Análise e projeto de circuitos eletrônicos utilizando softwares EDA. Análises experimentais de circuitos eletrônicos.'
$ws.Rows(15).RowHeight = 60

$ws.Range("A16").Value = 'Short syllabus:'
$ws.Range("B16").Value = 'Introduction to semiconductor materials and devices; Diode; Bipolar junction transistor; Field effect transistors; operational amplifiers; Small signal amplifiers; Power supplies. Analysis and design of electronic circuits using EDA software. Experimental analysis of electronic circuits.'
$ws.Range("C16").Value = 'Introduction to semiconductor materials and devices; Diode; Bipolar junction transistor; Field effect transistors; operational amplifiers; Small signal amplifiers; Power supplies. Analysis and design of electronic circuits using EDA software. Experimental analysis of electronic circuits.'
$ws.Rows(16).RowHeight = 60

$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = '1. Materiais Semicondutores; Diodos;    2. Retificadores de Tensão: Análise e Projeto;    3. Software EDA; Projeto de Placas de Circuito Impresso.    4. Transistor Bipolar de Junção (TBJ); Folha de Dados, Polarização e Chaveamento.    5. Transistor de Efeito de Campo (JFET - MOSFET); Folha de Dados, Polarização e Chaveamento.    6. Amplificadores para Pequenos Sinais: Análise e Projeto;    7. Amplificadores Operacionais: Buffer, Amplificação, Integração e Diferenciação;    8. Filtros Ativos: Análise e Projeto;    9. Fontes de Alimentação;    10. Circuitos Optoeletrônicos;'
$ws.Range("C17").Value = '1. Materiais Semicondutores; Diodos;    2. Retificadores de Tensão: Análise e Projeto;    3. Software EDA; Projeto de Placas de Circuito Impresso.    4. Transistor Bipolar de Junção (TBJ); Folha de Dados, Polarização e Chaveamento.    5. Transistor de Efeito de Campo (JFET - MOSFET); Folha de Dados, Polarização e Chaveamento.    6. Amplificadores para Pequenos Sinais: Análise e Projeto;    7. Amplificadores Operacionais: Buffer, Amplificação, Integração e Diferenciação;    8. Filtros Ativos: Análise e Projeto;    9. Fontes de Alimentação;    10. Circuitos Optoeletrônicos;'
$ws.Rows(17).RowHeight = 120

$ws.Range("A18").Value = 'Syllabus:'
$ws.Range("B18").Value = '1. Semiconductor Materials; Diodes;2. Voltage Rectifiers: Analysis and Design;3. EDA Software; Design of Printed Circuit Boards.4. Bipolar Junction Transistor (BJT); Data Sheet, Polarization and Switching.5. Field Effect Transistor (JFET - MOSFET); Data Sheet, Polarization and Switching.6. Small Signal Amplifiers: Analysis and Design;7. Operational Amplifiers: Buffer, Amplification, Integration and Differentiation;8. Active Filters: Analysis and Design;9. Power Supplies;10. Optoelectronic Circuits;'
$ws.Range("C18").Value = '1. Semiconductor Materials; Diodes;2. Voltage Rectifiers: Analysis and Design;3. EDA Software; Design of Printed Circuit Boards.4. Bipolar Junction Transistor (BJT); Data Sheet, Polarization and Switching.5. Field Effect Transistor (JFET - MOSFET); Data Sheet, Polarization and Switching.6. Small Signal Amplifiers: Analysis and Design;7. Operational Amplifiers: Buffer, Amplification, Integration and Differentiation;8. Active Filters: Analysis and Design;9. Power Supplies;10. Optoelectronic Circuits;'
$ws.Rows(18).RowHeight = 120

$ws.Range("A19").Value = 'Avaliação:'

$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas.'
$ws.Range("C20").Value = 'Aulas expositivas e práticas de laboratório com interações em grupo para a solução de problemas.'
$ws.Rows(20).RowHeight = 60

$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = 'Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T'
$ws.Range("C21").Value = 'Média aritmética (M) de provas individuais (P1 e P2) e trabalhos em grupo ao longo do semestre (T), tal que M = 0,3*P1+0,3*P2+0.4*T'
$ws.Rows(21).RowHeight = 60

$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C22").Value = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows(22).RowHeight = 60

$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = 'BROPHY, J. J. Eletrônica Básica. Guanabara Dois.NOVO, D. D. Eletrônica Aplicada. Editora da USP.SIMPSON, R.E. Introductory electronics for scientists and engineers. Allyn and Bacon.HOROWITZ, P.; HILL, W. The art of electronics. Cambridge University Press.MOTCHENBACHER, C. D.; FITCHEN, F.C. Low noise electronic design, John Wiley and Sons.MORRISON, R. Grounding and shielding techniques in instrumentation, John Wiley and Sons.ALEXANDER, C. K. E SADIKU, M. N. O. Fundamentos de Circuitos Elétricos. McGraw-Hill, 2013.NILSSON, J. W. E RIEDEL, S. A. Electric Circuits. Prentice Hall, 2011.BOYLESTAD, R. L. E NASHELSKY, L. Electronic Devices andCircuit Theory. Pearson, 2013'
$ws.Range("C23").Value = 'BROPHY, J. J. Eletrônica Básica. Guanabara Dois.NOVO, D. D. Eletrônica Aplicada. Editora da USP.SIMPSON, R.E. Introductory electronics for scientists and engineers. Allyn and Bacon.HOROWITZ, P.; HILL, W. The art of electronics. Cambridge University Press.MOTCHENBACHER, C. D.; FITCHEN, F.C. Low noise electronic design, John Wiley and Sons.MORRISON, R. Grounding and shielding techniques in instrumentation, John Wiley and Sons.ALEXANDER, C. K. E SADIKU, M. N. O. Fundamentos de Circuitos Elétricos. McGraw-Hill, 2013.NILSSON, J. W. E RIEDEL, S. A. Electric Circuits. Prentice Hall, 2011.BOYLESTAD, R. L. E NASHELSKY, L. Electronic Devices andCircuit Theory. Pearson, 2013'
$ws.Rows(23).RowHeight = 120

$ws.Range("A24").Value = 'Requisitos:'

$ws.Range("B25").Value = 'LOM3262 -  Circuitos Elétricos  (Requisito)
'
$ws.Range("C25").Value = 'LOM3262 -  Circuitos Elétricos  (Requisito)
'
$ws.Rows(25).RowHeight = 30

# --- Fix up column styles (A=style1, B=style2, C=style3) on the rebuilt rows ---
$ws.Range("AA1000").Copy()
$ws.Range("A10:A25").PasteSpecial(-4122)
$ws.Range("AA1001").Copy()
$ws.Range("B10:B25").PasteSpecial(-4122)
$ws.Range("AA1002").Copy()
$ws.Range("C10:C25").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Clean up the stashed templates ---
$ws.Range("AA1000:AA1002").Clear()

# --- Split the col-A / col-B width declaration (A stays 30.7109375, B keeps 60.7109375) ---
$wB = $ws.Columns.Item(2).ColumnWidth()
$ws.Columns.Item(2).ColumnWidth = $wB

